$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("F2").Value = 1.09
$ws.Range("J2").Value = 1.09
$ws.Range("N2").Value = 1.1
$ws.Range("P2").Value = 1.36
$ws.Range("S2").Value = 1.92
$ws.Range("V2").Value = 1.17
$ws.Range("W2").Value = 1.02

# Row 3
$ws.Range("F3").Value = 1.35
$ws.Range("J3").Value = 5.7
$ws.Range("N3").Value = 5.9
$ws.Range("P3").Value = 2.68
$ws.Range("Q3").Value = 1.56
$ws.Range("R3").Value = 1.67
$ws.Range("U3").Value = 2.04
$ws.Range("W3").Value = 3.75
$ws.Range("Y3").Value = 44
$ws.Range("AC3").Value = 13.5
$ws.Range("AD3").Value = 40
$ws.Range("AE3").Value = 160
$ws.Range("AH3").Value = 26

# Row 4
$ws.Range("G4").Value = 5.6
$ws.Range("N4").Value = 5.1
$ws.Range("P4").Value = 2.36
$ws.Range("W4").Value = 1.21
$ws.Range("AC4").Value = 9.199999999999999

# Row 5
$ws.Range("F5").Value = 1.04
$ws.Range("G5").Value = 1000
$ws.Range("H5").Value = 1.04
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 1.09
$ws.Range("L5").Value = 1.26
$ws.Range("N5").Value = 1.1
$ws.Range("P5").Value = 1.25
$ws.Range("Q5").Value = 1.32
$ws.Range("R5").Value = 1.16
$ws.Range("S5").Value = 2
$ws.Range("V5").Value = 1.01
$ws.Range("W5").Value = 1.01

# Row 6
$ws.Range("L6").Value = 1.26

# Row 8
$ws.Range("G8").Value = 8.4
$ws.Range("J8").Value = 5.4
$ws.Range("K8").Value = 5.6
$ws.Range("N8").Value = 6
$ws.Range("P8").Value = 2.74
$ws.Range("Q8").Value = 1.55
$ws.Range("S8").Value = 2.36
$ws.Range("T8").Value = 1.78
$ws.Range("U8").Value = 2.18
$ws.Range("V8").Value = 3.25
$ws.Range("W8").Value = 1.13
$ws.Range("X8").Value = 25
$ws.Range("Z8").Value = 10.5
$ws.Range("AN8").Value = 90

# Row 9
$ws.Range("F9").Value = 2.3
$ws.Range("G9").Value = 2.34
$ws.Range("I9").Value = 3.4
$ws.Range("O9").Value = 1.26
$ws.Range("P9").Value = 2.22
$ws.Range("Q9").Value = 1.79
$ws.Range("S9").Value = 2.96
$ws.Range("V9").Value = 1.41
$ws.Range("W9").Value = 1.74
$ws.Range("Y9").Value = 15.5
$ws.Range("AF9").Value = 16
$ws.Range("AH9").Value = 15

# Row 10
$ws.Range("F10").Value = 1.7
$ws.Range("P10").Value = 2.6

# Row 11
$ws.Range("H11").Value = 9.6
$ws.Range("I11").Value = 10
$ws.Range("J11").Value = 5.4
$ws.Range("K11").Value = 5.5
$ws.Range("O11").Value = 1.18
$ws.Range("R11").Value = 1.69
$ws.Range("S11").Value = 2.4
$ws.Range("W11").Value = 3.45
$ws.Range("Z11").Value = 90
$ws.Range("AD11").Value = 34

# Row 12
$ws.Range("F12").Value = 1.34
$ws.Range("G12").Value = 1.35
$ws.Range("H12").Value = 9.6
$ws.Range("K12").Value = 6.8
$ws.Range("T12").Value = 1.66
$ws.Range("AJ12").Value = 13
$ws.Range("AN12").Value = 3.4

# Row 13
$ws.Range("P13").Value = 2.06
$ws.Range("Q13").Value = 1.91
$ws.Range("T13").Value = 1.72
$ws.Range("U13").Value = 2.3
$ws.Range("AI13").Value = 44
$ws.Range("AL13").Value = 36
$ws.Range("AM13").Value = 80

# Row 14
$ws.Range("F14").Value = 2
$ws.Range("L14").Value = 1.39
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 4
$ws.Range("O14").Value = 1.28
$ws.Range("R14").Value = 1.36
$ws.Range("S14").Value = 2.8
$ws.Range("T14").Value = 1.6
$ws.Range("U14").Value = 2.04
$ws.Range("V14").Value = 1.29
$ws.Range("W14").Value = 1.86
$ws.Range("X14").Value = 1000
$ws.Range("Y14").Value = 22
$ws.Range("Z14").Value = 42
$ws.Range("AA14").Value = 100
$ws.Range("AB14").Value = 15
$ws.Range("AC14").Value = 12
$ws.Range("AD14").Value = 22
$ws.Range("AE14").Value = 60
$ws.Range("AF14").Value = 19.5
$ws.Range("AG14").Value = 15
$ws.Range("AH14").Value = 24
$ws.Range("AI14").Value = 70
$ws.Range("AJ14").Value = 34
$ws.Range("AK14").Value = 30
$ws.Range("AL14").Value = 46
$ws.Range("AM14").Value = 1000
$ws.Range("AN14").Value = 19
$ws.Range("AO14").Value = 55

# Row 15
$ws.Range("F15").Value = 1.68
$ws.Range("G15").Value = 2.1
$ws.Range("H15").Value = 4.4
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 3.2
$ws.Range("K15").Value = 5.7
$ws.Range("P15").Value = 1.55
$ws.Range("Q15").Value = 2.12
